$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 405, shifting existing rows 405:486 down to 406:487
$ws.Rows.Item(405).Insert()

# Populate the new row 405 with the new data record
$ws.Cells.Item(405, 1).Value = 4
$ws.Cells.Item(405, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(405, 3).Value = "Los Lagos"
$ws.Cells.Item(405, 4).Value = 44641
$ws.Cells.Item(405, 5).Value = 10
$ws.Cells.Item(405, 6).Value = 100112004
$ws.Cells.Item(405, 7).Value = "Cebolla"
$ws.Cells.Item(405, 8).Value = "Sin especificar"
$ws.Cells.Item(405, 9).Value = "1a (cosecha)"
$ws.Cells.Item(405, 10).Value = 250
$ws.Cells.Item(405, 11).Value = 10000
$ws.Cells.Item(405, 12).Value = 10000
$ws.Cells.Item(405, 13).Value = 10000
$ws.Cells.Item(405, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(405, 15).Value = "Región del Maule"
$ws.Cells.Item(405, 16).Value = 400
$ws.Cells.Item(405, 17).Value = 25
$ws.Cells.Item(405, 18).Value = "Hortaliza"

# Match the date-style formatting used by column D (numFmtId 165) on the new row
$ws.Cells.Item(405, 4).NumberFormat = $ws.Cells.Item(406, 4).NumberFormat
